$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.252965807914734
$ws.Range("B1").Value = 2.364249706268311
$ws.Range("C1").Value = 3.366795778274536
$ws.Range("D1").Value = 2.46784234046936
$ws.Range("E1").Value = 1.361463665962219
